$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Authentication")
$ws.Range("A1").Value = "Hello"
